# "updatedd the simple data generation procedure"
# Change the simple-data generation parameter "allocation_min_qty_supplier"
# (cell B4 on Sheet1) from 0.2 to 0.1, and leave the selection on that cell
# (matching Excel's own behaviour of remembering the last-edited cell as the
# active selection when the workbook is next saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 0.1
$ws.Range("B4").Select()
